$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - header row
$ws.Range("A1").Value = "Thời gian"
$ws.Range("B1").Value = "Thứ Hai"
$ws.Range("C1").Value = "Thứ Ba"
$ws.Range("D1").Value = "Thứ Tư"
$ws.Range("E1").Value = "Thứ Năm"
$ws.Range("F1").Value = "Thứ Sáu"

# G1 ("Chủ Nhật") shifts right to H1; G1 becomes the new "Thứ Bảy" column.
# Clone G1's header formatting onto H1 before overwriting either value.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Thứ Bảy"
$ws.Range("H1").Value = "Chủ Nhật"

# Row 2
$ws.Range("A2").Value = "06:00 – 08:00"
$ws.Range("B2").Value = "Lập trình hướng đối tượng"
$ws.Range("C2").Value = "Trí tuệ nhân tạo"
$ws.Range("D2").Value = "Nguyên lý Hệ điều hành"
$ws.Range("E2").Value = "Tư tưởng HCM"
$ws.Range("F2").Value = "Tiếng Anh 1"
$ws.Range("G2").Value = "(Đang đi làm)"
$ws.Range("H2").Value = "Ôn tập nhanh tất cả môn (2h)"

# Row 3
$ws.Range("A3").Value = "08:30 – 10:30"
$ws.Range("B3").Value = "Trí tuệ nhân tạo"
$ws.Range("C3").Value = "Nguyên lý HĐH"
$ws.Range("D3").Value = "Tư tưởng HCM"
$ws.Range("E3").Value = "Tiếng Anh 1"
$ws.Range("F3").Value = "Lập trình HĐTĐ"
$ws.Range("G3").Value = "(Đang đi làm)"
$ws.Range("H3").Value = "Hoàn thiện dự án nhỏ hàng tuần (2h)"

# Row 4
$ws.Range("A4").Value = "14:00 – 16:00"
$ws.Range("B4").Value = "Nguyên lý Hệ điều hành"
$ws.Range("C4").Value = "Tư tưởng HCM"
$ws.Range("D4").Value = "Tiếng Anh 1"
$ws.Range("E4").Value = "Lập trình HĐTĐ"
$ws.Range("F4").Value = "Trí tuệ nhân tạo"
$ws.Range("G4").Value = "(Đang đi làm đến 17:00)"
$ws.Range("H4").Value = "Buffer & ôn chuyên sâu (2h)"

# Row 5
$ws.Range("A5").Value = "16:30 – 18:00"
$ws.Range("B5").Value = "Tư tưởng HCM"
$ws.Range("C5").Value = "Tiếng Anh 1"
$ws.Range("D5").Value = "Lập trình HĐTĐ"
$ws.Range("E5").Value = "Trí tuệ nhân tạo"
$ws.Range("F5").Value = "Nguyên lý HĐH"
$ws.Range("G5").Value = "17:00 – 18:30: Gym18:30 – 20:00: Ăn tối/ nghỉ"
$ws.Range("H5").Value = "Nghỉ/nghỉ linh hoạt"

# Row 6
$ws.Range("A6").Value = "20:00 – 22:00"
$ws.Range("B6").Value = "Gym (tối)"
$ws.Range("C6").Value = "Gym (tối)"
$ws.Range("D6").Value = "Gym (tối)"
$ws.Range("E6").Value = "Gym (tối)"
$ws.Range("F6").Value = "Gym (tối)"
$ws.Range("G6").Value = "20:00 – 22:00: Dự án nhỏ hàng tuần"
$ws.Range("H6").Value = "20:00 – 21:30: Gym (tuỳ chọn)"

# Remove old leftover rows (7-11), which are no longer part of the table
$ws.Range("A7:G11").Delete()
